$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "FaultInjectionLidar1"
$ws.Range("B21").Value = 0
$ws.Range("A22").Value = "FaultInjectionLidar2"
$ws.Range("B22").Value = 0
$ws.Range("A23").Value = "FaultInjectionLidar3"
$ws.Range("B23").Value = 0
$ws.Range("A24").Value = "FaultInjectionLidar4"
$ws.Range("B24").Value = 0

$ws.Range("A21:B24").Select()
